# Updated some column header definitions
# - Delete the StationData sheet entirely.
# - Update several attribute definitions on the ColumnHeaders sheet
#   (bottle_other_method, alternate_sample_id, filter_size, chl, phaeo,
#   nearest_station, quality_flag).
# - Add a role for Pierre Marrec on the Personnel sheet.
# - Leave Personnel as the active sheet / tab (matches final saved state).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- ColumnHeaders sheet: refresh a handful of attribute definitions ---
$colHeaders = $wb.Worksheets.Item("ColumnHeaders")

$colHeaders.Range("B5").Value = "Associated niskin with samples analyzed by another method"
$colHeaders.Range("B12").Value = "Alternate identifier for the same sample"
$colHeaders.Range("B13").Value = "Size fraction with filter and when applicable prefilter in micrometers "
$colHeaders.Range("B14").Value = "Concentration of chlorophyll-a per unit volume of the water body http://vocab.nerc.ac.uk/collection/P01/current/CPHLZZXX/ corresponding to BCO-DMO chl_a SeaBASS Chl and CF mass_concentration_of_ chlorophyll_a_in_sea_water"
$colHeaders.Range("B15").Value = "Concentration of phaeopigments per unit volume of the water body"
$colHeaders.Range("B17").Value = "NES-LTER standard station nearest to the sample location"
$colHeaders.Range("B19").Value = "Quality flag provided for a subset of the data"

# --- Personnel sheet: Pierre Marrec's role ---
$personnel = $wb.Worksheets.Item("Personnel")
$personnel.Range("G9").Value = "postdoctoralResearcher"

# --- Remove the StationData sheet (its rows were folded into ColumnHeaders) ---
$stationData = $wb.Worksheets.Item("StationData")
[void]$stationData.Delete()

# --- View state: ColumnHeaders deselected, Personnel becomes the active tab ---
[void]$colHeaders.Activate()
$colHeaders.Range("B19").Select() | Out-Null

[void]$personnel.Activate()
$personnel.Range("G13").Select() | Out-Null
